$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 4: description in column A, URL (as a hyperlink) in column B.
# Write the URL cell before the description cell so the shared-strings
# table picks up the same ordering as the authored workbook (URL first).
$url = "https://stackoverflow.com/questions/35655361/angular2-how-to-load-data-before-rendering-the-component"
$desc = "Getting data syncronously / Getting data after page loads so  page doesn't show data"

$ws.Range("B4").Value2 = $url
$ws.Range("A4").Value2 = $desc

# Turn B4 into a real hyperlink pointing at the URL it displays.
$ws.Hyperlinks.Add($ws.Range("B4"), $url) | Out-Null

# Hyperlinks.Add mints its own cell format; put B4 back on the shared
# "Hyperlink" cell style so it matches the other link cells in the column.
$ws.Range("B4").Style = "Hyperlink"

# Match the author's final selection (cell A4).
$ws.Range("A4").Select()
